$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27; this shifts rows 27-93 down to 28-94.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly record.
$ws.Cells.Item(27, 1).Value = 4
$ws.Cells.Item(27, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(27, 3).Value = "Los Lagos"
$ws.Cells.Item(27, 4).Value = 44498
$ws.Cells.Item(27, 5).Value = 10
$ws.Cells.Item(27, 6).Value = "Fruta"
$ws.Cells.Item(27, 7).Value = 100108
$ws.Cells.Item(27, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(27, 9).Value = 100108002
$ws.Cells.Item(27, 10).Value = "Mango"
$ws.Cells.Item(27, 11).Value = "Sin especificar"
$ws.Cells.Item(27, 12).Value = "Primera"
$ws.Cells.Item(27, 13).Value = 300
$ws.Cells.Item(27, 14).Value = 8000
$ws.Cells.Item(27, 15).Value = 8500
$ws.Cells.Item(27, 16).Value = 8250
$ws.Cells.Item(27, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(27, 18).Value = "Perú"
$ws.Cells.Item(27, 19).Value = 2062
$ws.Cells.Item(27, 20).Value = 4

# Apply the same date cell style (format) as the rest of column D to the new D27 cell.
$ws.Cells.Item(27, 4).NumberFormat = $ws.Cells.Item(28, 4).NumberFormat
